$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths ---
# Column B (Thema/Aufgabe text) becomes much wider; columns C:G get a touch narrower.
$ws.Columns("B").ColumnWidth = 53.944010416666664
$ws.Columns("C:G").ColumnWidth = 12.721354166666666

# --- Row 10: new task "Recherge/Implementieren des adaptiven Headers" ---
$ws.Range("B10").Value = "Recherge/Implementieren des adaptiven Headers"
$ws.Range("C10").Value = "x"
$ws.Range("D10").Clear()
$ws.Range("E10").Value = 43852
$ws.Range("E10").NumberFormat = "m/d/yy"
$ws.Range("F10").Value = 43858
$ws.Range("F10").NumberFormat = "m/d/yy"

# --- Row 11: new task "Navbar implementiren + Design erweitern" ---
$ws.Range("B11").Value = "Navbar implementiren + Design erweitern"
$ws.Range("D11").Value = "x"
$ws.Range("E11").Value = 43852
$ws.Range("E11").NumberFormat = "m/d/yy"
$ws.Range("F11").Value = 43858
$ws.Range("F11").NumberFormat = "m/d/yy"
